$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 116
$ws.Range("H116").Value = 3353.077
$ws.Range("I116").Value = 2265
$ws.Range("J116").Value = 4622.5
$ws.Range("K116").Value = 2265
$ws.Range("L116").Value = 4622.5
$ws.Range("M116").Value = 1177
$ws.Range("N116").Value = -11506.5
# Row 125
$ws.Range("H125").Value = 1697
$ws.Range("I125").Value = 1133
$ws.Range("J125").Value = 2120
$ws.Range("K125").Value = 10197
$ws.Range("L125").Value = 19080
$ws.Range("M125").Value = -7737
$ws.Range("N125").Value = -24000
# Row 129
$ws.Range("H129").Value = 1449.2084
$ws.Range("I129").Value = 813
$ws.Range("J129").Value = 1711.1765
$ws.Range("K129").Value = 2439
$ws.Range("L129").Value = 5133.529500000001
$ws.Range("M129").Value = 2561
$ws.Range("N129").Value = -15133.5295
# Row 137
$ws.Range("H137").Value = 2032.3846
$ws.Range("I137").Value = 2099.4348
$ws.Range("J137").Value = 1936
$ws.Range("K137").Value = 6298.3044
$ws.Range("L137").Value = 5808
$ws.Range("M137").Value = -3748.3044
$ws.Range("N137").Value = -10908
# Row 138
$ws.Range("H138").Value = 3210.279
$ws.Range("I138").Value = 1566.6538
$ws.Range("J138").Value = 5724.0586
$ws.Range("K138").Value = 4699.9614
$ws.Range("L138").Value = 17172.1758
$ws.Range("M138").Value = 440.0385999999999
$ws.Range("N138").Value = -27452.1758
# Row 141
$ws.Range("H141").Value = 14147.158
$ws.Range("I141").Value = 1891.2307
$ws.Range("J141").Value = 40701.668
$ws.Range("K141").Value = 5673.6921
$ws.Range("L141").Value = 122105.004
$ws.Range("M141").Value = -493.6921000000002
$ws.Range("N141").Value = -132465.004

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4846.242
$ws.Range("I32").Value = 4478.011
$ws.Range("K32").Value = 4478.011
$ws.Range("M32").Value = -4191.011

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
# Row 41
$ws.Range("H41").Value = 229386
$ws.Range("J41").Value = 229386
$ws.Range("L41").Value = 229386
$ws.Range("N41").Value = -230162
# Row 134
$ws.Range("H134").Value = 4618.5283
$ws.Range("I134").Value = 1902.5518
$ws.Range("J134").Value = 7900.3335
$ws.Range("K134").Value = 5707.6554
$ws.Range("L134").Value = 23701.0005
$ws.Range("M134").Value = -3172.6554
$ws.Range("N134").Value = -28771.0005

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7938745
$ws.Range("I31").Value = 1560.08
$ws.Range("J31").Value = 19611076
$ws.Range("K31").Value = 1560.08
$ws.Range("L31").Value = 19611076
$ws.Range("M31").Value = -1265.08
$ws.Range("N31").Value = -19611666
# Row 34
$ws.Range("H34").Value = 7938745
$ws.Range("I34").Value = 1560.08
$ws.Range("J34").Value = 19611076
$ws.Range("K34").Value = 1560.08
$ws.Range("L34").Value = 19611076
$ws.Range("M34").Value = -1358.08
$ws.Range("N34").Value = -19611480
# Row 132
$ws.Range("H132").Value = 2686.7778
$ws.Range("I132").Value = 1772.8276
$ws.Range("J132").Value = 6473.143
$ws.Range("K132").Value = 5318.4828
$ws.Range("L132").Value = 19419.429
$ws.Range("M132").Value = -2788.4828
$ws.Range("N132").Value = -24479.429

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 62
$ws.Range("H62").Value = 2720.1
$ws.Range("I62").Value = 999.5
$ws.Range("J62").Value = 3150.25
$ws.Range("K62").Value = 2998.5
$ws.Range("L62").Value = 9450.75
$ws.Range("M62").Value = -2312.5
$ws.Range("N62").Value = -10822.75
# Row 63
$ws.Range("H63").Value = 2668.6667
$ws.Range("I63").Value = 870.6667
$ws.Range("J63").Value = 4466.6665
$ws.Range("K63").Value = 2612.0001
$ws.Range("L63").Value = 13399.9995
$ws.Range("M63").Value = -1863.0001
$ws.Range("N63").Value = -14897.9995
# Row 65
$ws.Range("H65").Value = 2720.1
$ws.Range("I65").Value = 999.5
$ws.Range("J65").Value = 3150.25
$ws.Range("K65").Value = 8995.5
$ws.Range("L65").Value = 28352.25
$ws.Range("M65").Value = -5563.5
$ws.Range("N65").Value = -35216.25
# Row 66
$ws.Range("H66").Value = 2668.6667
$ws.Range("I66").Value = 870.6667
$ws.Range("J66").Value = 4466.6665
$ws.Range("K66").Value = 7836.0003
$ws.Range("L66").Value = 40199.9985
$ws.Range("M66").Value = -4092.0003
$ws.Range("N66").Value = -47687.9985

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 75
$ws.Range("H75").Value = 20001
$ws.Range("J75").Value = 20001
$ws.Range("L75").Value = 20001
$ws.Range("N75").Value = -21749
# Row 78
$ws.Range("H78").Value = 20001
$ws.Range("J78").Value = 20001
$ws.Range("L78").Value = 60003
$ws.Range("N78").Value = -68739
# Row 86
$ws.Range("H86").Value = 17024.334
$ws.Range("J86").Value = 17024.334
$ws.Range("L86").Value = 17024.334
$ws.Range("N86").Value = -19396.334
# Row 89
$ws.Range("H89").Value = 17024.334
$ws.Range("J89").Value = 17024.334
$ws.Range("L89").Value = 51073.00199999999
$ws.Range("N89").Value = -62929.00199999999
# Row 107
$ws.Range("H107").Value = 492
$ws.Range("I107").Value = 478.29413
$ws.Range("J107").Value = 530.8333
$ws.Range("K107").Value = 478.29413
$ws.Range("L107").Value = 530.8333
$ws.Range("M107").Value = 1441.70587
$ws.Range("N107").Value = -4370.8333
# Row 132
$ws.Range("H132").Value = 1159983.4
$ws.Range("I132").Value = 2606119
$ws.Range("J132").Value = 3074.9
$ws.Range("K132").Value = 7818357
$ws.Range("L132").Value = 9224.700000000001
$ws.Range("M132").Value = -7815827
$ws.Range("N132").Value = -14284.7

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1741.25
$ws.Range("I22").Value = 1790
$ws.Range("J22").Value = 1692.5
$ws.Range("K22").Value = 1790
$ws.Range("L22").Value = 1692.5
$ws.Range("M22").Value = -1495
$ws.Range("N22").Value = -2282.5
# Row 27
$ws.Range("H27").Value = 1741.25
$ws.Range("I27").Value = 1790
$ws.Range("J27").Value = 1692.5
$ws.Range("K27").Value = 1790
$ws.Range("L27").Value = 1692.5
$ws.Range("M27").Value = -1683
$ws.Range("N27").Value = -1906.5
# Row 40
$ws.Range("H40").Value = 3520.7307
$ws.Range("I40").Value = 2651.3333
$ws.Range("J40").Value = 5476.875
$ws.Range("K40").Value = 2651.3333
$ws.Range("L40").Value = 5476.875
$ws.Range("M40").Value = -2515.3333
$ws.Range("N40").Value = -5748.875
# Row 46
$ws.Range("H46").Value = 2001637.8
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2001637.8
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2001637.8
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -2002013.8
# Row 101
$ws.Range("H101").Value = 12792.625
$ws.Range("J101").Value = 12792.625
$ws.Range("L101").Value = 12792.625
$ws.Range("N101").Value = -19282.625
# Row 122
$ws.Range("H122").Value = 3876.4443
$ws.Range("I122").Value = 3857.6
$ws.Range("J122").Value = 3900
$ws.Range("K122").Value = 11572.8
$ws.Range("L122").Value = 11700
$ws.Range("M122").Value = -9122.799999999999
$ws.Range("N122").Value = -16600

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 3950
$ws.Range("I39").Value = 2900
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 2900
$ws.Range("L39").Value = 5000
$ws.Range("M39").Value = -2487
$ws.Range("N39").Value = -5826
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
# Row 52
$ws.Range("H52").Value = 6000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 6000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 6000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -6452
# Row 54
$ws.Range("H54").Value = 9399.799999999999
$ws.Range("J54").Value = 10999.75
$ws.Range("L54").Value = 10999.75
$ws.Range("N54").Value = -12039.75
# Row 81
$ws.Range("H81").Value = 2216
$ws.Range("I81").Value = 2427.4285
$ws.Range("J81").Value = 1722.6666
$ws.Range("K81").Value = 4854.857
$ws.Range("L81").Value = 3445.3332
$ws.Range("M81").Value = -3793.857
$ws.Range("N81").Value = -5567.3332
# Row 84
$ws.Range("H84").Value = 2216
$ws.Range("I84").Value = 2427.4285
$ws.Range("J84").Value = 1722.6666
$ws.Range("K84").Value = 24274.285
$ws.Range("L84").Value = 17226.666
$ws.Range("M84").Value = -18970.285
$ws.Range("N84").Value = -27834.666
# Row 103
$ws.Range("H103").Value = 20696
$ws.Range("J103").Value = 20696
$ws.Range("L103").Value = 20696
$ws.Range("N103").Value = -23040
# Row 107
$ws.Range("H107").Value = 539
$ws.Range("I107").Value = 516.125
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 1548.375
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 371.625
$ws.Range("N107").Value = -5640
# Row 113
$ws.Range("H113").Value = 929.2222
$ws.Range("I113").Value = 486.84616
$ws.Range("J113").Value = 1340
$ws.Range("K113").Value = 1460.53848
$ws.Range("L113").Value = 4020
$ws.Range("M113").Value = 709.4615200000001
$ws.Range("N113").Value = -8360
# Row 132
$ws.Range("H132").Value = 1885.6182
$ws.Range("I132").Value = 1493.826
$ws.Range("J132").Value = 3888.111
$ws.Range("K132").Value = 4481.478
$ws.Range("L132").Value = 11664.333
$ws.Range("M132").Value = -1951.478
$ws.Range("N132").Value = -16724.333
